# Update Leve profit-tracking figures across all item sheets
# (scheduled runner refresh of currentAveragePrice / LevePrice* / LeveProfit* columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2733373.8
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2733373.8
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 8200121.399999999
$ws.Range("N17").Value = -8200457.399999999
$ws.Range("M17").ClearContents()

$ws.Range("H86").Value = 9285.571
$ws.Range("I86").Value = 6899.8
$ws.Range("K86").Value = 6899.8
$ws.Range("M86").Value = -5776.8

$ws.Range("H88").Value = 5298.625
$ws.Range("I88").Value = 7759.8
$ws.Range("J88").Value = 4179.909
$ws.Range("K88").Value = 7759.8
$ws.Range("L88").Value = 4179.909
$ws.Range("M88").Value = -7353.8
$ws.Range("N88").Value = -4991.909

$ws.Range("H89").Value = 9285.571
$ws.Range("I89").Value = 6899.8
$ws.Range("K89").Value = 34499
$ws.Range("M89").Value = -28883

$ws.Range("H91").Value = 5298.625
$ws.Range("I91").Value = 7759.8
$ws.Range("J91").Value = 4179.909
$ws.Range("K91").Value = 7759.8
$ws.Range("L91").Value = 4179.909
$ws.Range("M91").Value = -6355.8
$ws.Range("N91").Value = -6987.909

$ws.Range("H133").Value = 15199.8
$ws.Range("J133").Value = 15199.8
$ws.Range("L133").Value = 15199.8
$ws.Range("N133").Value = -25319.8

$ws.Range("H136").Value = 68980.586
$ws.Range("J136").Value = 68980.586
$ws.Range("L136").Value = 68980.586
$ws.Range("N136").Value = -79180.586

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3419.4075
$ws.Range("I45").Value = 3568.2
$ws.Range("J45").Value = 2994.2856
$ws.Range("K45").Value = 3568.2
$ws.Range("L45").Value = 2994.2856
$ws.Range("M45").Value = -3191.2
$ws.Range("N45").Value = -3748.2856

$ws.Range("H97").Value = 744.41174
$ws.Range("I97").Value = 781.1539
$ws.Range("K97").Value = 781.1539
$ws.Range("M97").Value = -285.1539

$ws.Range("H132").Value = 15325.333
$ws.Range("I132").Value = 3882.9
$ws.Range("K132").Value = 11648.7
$ws.Range("M132").Value = -9118.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2217.7334
$ws.Range("I22").Value = 2365.9285
$ws.Range("J22").Value = 143
$ws.Range("K22").Value = 2365.9285
$ws.Range("L22").Value = 143
$ws.Range("M22").Value = -2192.9285
$ws.Range("N22").Value = -489

$ws.Range("H94").Value = 1252.5264
$ws.Range("I94").Value = 788.975
$ws.Range("J94").Value = 2343.2354
$ws.Range("K94").Value = 788.975
$ws.Range("L94").Value = 2343.2354
$ws.Range("M94").Value = -337.975
$ws.Range("N94").Value = -3245.2354

$ws.Range("H99").Value = 27507.95
$ws.Range("J99").Value = 16892.715
$ws.Range("L99").Value = 16892.715
$ws.Range("N99").Value = -19888.715

$ws.Range("H134").Value = 31316.809
$ws.Range("I134").Value = 33630.156
$ws.Range("J134").Value = 23914.1
$ws.Range("K134").Value = 100890.468
$ws.Range("L134").Value = 71742.29999999999
$ws.Range("M134").Value = -98355.46800000001
$ws.Range("N134").Value = -76812.29999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30062
$ws.Range("I31").Value = 14833.75
$ws.Range("J31").Value = 42244.6
$ws.Range("K31").Value = 14833.75
$ws.Range("L31").Value = 42244.6
$ws.Range("M31").Value = -14538.75
$ws.Range("N31").Value = -42834.6

$ws.Range("H34").Value = 30062
$ws.Range("I34").Value = 14833.75
$ws.Range("J34").Value = 42244.6
$ws.Range("K34").Value = 14833.75
$ws.Range("L34").Value = 42244.6
$ws.Range("M34").Value = -14631.75
$ws.Range("N34").Value = -42648.6

$ws.Range("H80").Value = 14000
$ws.Range("J80").Value = 14000
$ws.Range("L80").Value = 14000
$ws.Range("N80").Value = -16246

$ws.Range("H83").Value = 14000
$ws.Range("J83").Value = 14000
$ws.Range("L83").Value = 42000
$ws.Range("N83").Value = -53232

$ws.Range("H94").Value = 5213.3076
$ws.Range("J94").Value = 655.2727
$ws.Range("L94").Value = 655.2727
$ws.Range("N94").Value = -1557.2727

$ws.Range("H99").Value = 339612.72
$ws.Range("I99").Value = 486010.5
$ws.Range("K99").Value = 486010.5
$ws.Range("M99").Value = -484512.5

$ws.Range("H126").Value = 339612.72
$ws.Range("I126").Value = 486010.5
$ws.Range("K126").Value = 1458031.5
$ws.Range("M126").Value = -1455561.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1957.6666
$ws.Range("I68").Value = 863.3333
$ws.Range("J68").Value = 2322.4443
$ws.Range("K68").Value = 2589.9999
$ws.Range("L68").Value = 6967.3329
$ws.Range("M68").Value = -1778.9999
$ws.Range("N68").Value = -8589.332900000001

$ws.Range("H71").Value = 1957.6666
$ws.Range("I71").Value = 863.3333
$ws.Range("J71").Value = 2322.4443
$ws.Range("K71").Value = 7769.9997
$ws.Range("L71").Value = 20901.9987
$ws.Range("M71").Value = -3713.9997
$ws.Range("N71").Value = -29013.9987

$ws.Range("H137").Value = 5097.4
$ws.Range("J137").Value = 4664
$ws.Range("L137").Value = 13992
$ws.Range("N137").Value = -24192

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3824.6875
$ws.Range("I102").Value = 4432.9585
$ws.Range("J102").Value = 1999.875
$ws.Range("K102").Value = 4432.9585
$ws.Range("L102").Value = 1999.875
$ws.Range("M102").Value = -2810.9585
$ws.Range("N102").Value = -5243.875

$ws.Range("H126").Value = 11938.1875
$ws.Range("I126").Value = 19710.334
$ws.Range("K126").Value = 59131.00199999999
$ws.Range("M126").Value = -56661.00199999999

$ws.Range("H132").Value = 11657.429
$ws.Range("I132").Value = 9459.593000000001
$ws.Range("J132").Value = 19075.125
$ws.Range("K132").Value = 28378.779
$ws.Range("L132").Value = 57225.375
$ws.Range("M132").Value = -25848.779
$ws.Range("N132").Value = -62285.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2547.3513
$ws.Range("I22").Value = 1941.6666
$ws.Range("J22").Value = 3342.3125
$ws.Range("K22").Value = 1941.6666
$ws.Range("L22").Value = 3342.3125
$ws.Range("M22").Value = -1646.6666
$ws.Range("N22").Value = -3932.3125

$ws.Range("H27").Value = 2547.3513
$ws.Range("I27").Value = 1941.6666
$ws.Range("J27").Value = 3342.3125
$ws.Range("K27").Value = 1941.6666
$ws.Range("L27").Value = 3342.3125
$ws.Range("M27").Value = -1834.6666
$ws.Range("N27").Value = -3556.3125

$ws.Range("H40").Value = 8486.333000000001
$ws.Range("I40").Value = 2500
$ws.Range("K40").Value = 2500
$ws.Range("M40").Value = -2364

$ws.Range("H61").Value = 4061.3572
$ws.Range("I61").Value = 3068.75
$ws.Range("K61").Value = 3068.75
$ws.Range("M61").Value = -2866.75

$ws.Range("H68").Value = 2538.6
$ws.Range("I68").Value = 2031.6666
$ws.Range("K68").Value = 2031.6666
$ws.Range("M68").Value = -1282.6666

$ws.Range("H71").Value = 2538.6
$ws.Range("I71").Value = 2031.6666
$ws.Range("K71").Value = 10158.333
$ws.Range("M71").Value = -6414.333000000001

$ws.Range("H93").Value = 5440.0645
$ws.Range("I93").Value = 4023.762
$ws.Range("J93").Value = 8414.299999999999
$ws.Range("K93").Value = 4023.762
$ws.Range("L93").Value = 8414.299999999999
$ws.Range("M93").Value = -2775.762
$ws.Range("N93").Value = -10910.3

$ws.Range("H113").Value = 4061.3572
$ws.Range("I113").Value = 3068.75
$ws.Range("K113").Value = 3068.75
$ws.Range("M113").Value = -898.75

$ws.Range("H122").Value = 6425.2974
$ws.Range("J122").Value = 6927.0557
$ws.Range("L122").Value = 20781.1671
$ws.Range("N122").Value = -25681.1671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

$ws.Range("H107").Value = 1306.8572
$ws.Range("I107").Value = 1556.2222
$ws.Range("J107").Value = 858
$ws.Range("K107").Value = 4668.6666
$ws.Range("L107").Value = 2574
$ws.Range("M107").Value = -2748.6666
$ws.Range("N107").Value = -6414

$ws.Range("H118").Value = 78999.5
$ws.Range("J118").Value = 78999.5
$ws.Range("L118").Value = 78999.5
$ws.Range("N118").Value = -82313.5

$ws.Range("H122").Value = 3714.5806
$ws.Range("I122").Value = 2811.6
$ws.Range("K122").Value = 8434.799999999999
$ws.Range("M122").Value = -5984.799999999999
